$d = $word.ActiveDocument

# The assignment's due date changed from March 6th to March 9th.
#
# NB: this runtime's Find/Replace (wdReplaceAll / wdReplaceOne) searches the
# whole document story regardless of which Range it was invoked on, so a
# plain text search-and-replace for "6th" -> "9h" would also clobber the
# (extremely common) substring "th" everywhere else in the document. Instead,
# locate the "Due:" paragraph once (read-only Find, no replacement - safe),
# then edit the two affected characters directly via explicit Range offsets
# so only that paragraph is touched.

$seek = $d.Content
$found = $seek.Find.Execute("Due:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $dueParagraph = $seek.Paragraphs(1)
    $pStart = $dueParagraph.Range.Start
    $pText = $dueParagraph.Range.Text

    $idx = $pText.IndexOf("6th")
    if ($idx -ge 0) {
        # "6" -> "9" (keeps its own run / formatting)
        $dayRange = $d.Range($pStart + $idx, $pStart + $idx + 1)
        $dayRange.Text = "9"

        # "th" -> "h" (the ordinal suffix run, still superscript)
        $suffixRange = $d.Range($pStart + $idx + 1, $pStart + $idx + 3)
        $suffixRange.Text = "h"

        Write-Host "Due paragraph now:" $dueParagraph.Range.Text
    }
}
